# Revised function for determining ERC vs. OT win in # of grammar changes,
# to accurately account for ties (equal counts -> "N/A" instead of "ERC").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nonnoisy_changeinfo_original")

for ($r = 2; $r -le 67; $r++) {
    $formula = '=_xlfn.IFS(AND(B' + $r + '<C' + $r + ', AND(NOT(ISTEXT(B' + $r + ')), NOT(ISTEXT(C' + $r + ')))),"OT", AND(B' + $r + '>C' + $r + ', AND(NOT(ISTEXT(B' + $r + ')), NOT(ISTEXT(C' + $r + ')))), "ERC", TRUE, "N/A")'
    $ws.Cells.Item($r, 4).FormulaArray = $formula
}

# Recalculate so the new formula results are live before re-sorting the table.
$excel.Calculate()

# Re-sort the data (excluding the header row) by the outcome column so the
# table stays ordered ERC / N-A / OT, now that ties land in "N/A".
$dataRange = $ws.Range("A2:D67")
$dataRange.Sort($ws.Range("D2")) | Out-Null

# The re-save now leaves sheet1 as the active/selected sheet.
$ws.Select() | Out-Null
$ws.Range("F46").Select() | Out-Null
